# DaySale report update:
#  - a new item "THYROXINE 50MCG 100 TAB." is inserted into the shortage
#    list (alphabetically it lands right before the row that used to be
#    item #15, "سرنجات 10 سم"), so every item from the old #15 onward
#    shifts down by one row.
#  - the totals row (سعر البيع column) grows by the new item's price.
#  - the footer timestamp is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Force a cell to be written as literal text even when its number
    # format would otherwise make Excel coerce the string into a number
    # (several columns in this report store numeric-looking data as text).
    param($range, [string]$text)
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.NumberFormat = $fmt
}

# ------------------------------------------------------------------
# 1. Make room: the workbook currently ends with
#       row 26 = totals ("سعر البيع" sum)
#       row 27 = footer (timestamp / page / developer credit)
#    After the edit there is one more item row, so the totals row becomes
#    row 27 and the footer becomes row 28. Insert a blank row at 27 so the
#    existing totals row (26) and footer row (27->28) keep their own
#    formatting untouched, and only the (currently blank) new row needs
#    styling.
# ------------------------------------------------------------------
$ws.Rows("27").Insert()
$ws.Rows("27").RowHeight = 25.5

# ------------------------------------------------------------------
# 2. Shift the last five existing items (old rows 21-25) down one row
#    (new rows 22-26), carrying their text forward. Work bottom-up so we
#    never overwrite a row before reading it. Row 26 first needs the
#    normal item-row formatting (it is still styled as the totals row at
#    this point), which we clone from row 25. Q26 is still part of the
#    (soon to be retired) P26:Q26 totals merge, and Excel refuses writes
#    to the non-anchor cell of a merged range, so unmerge it first.
# ------------------------------------------------------------------
$ws.Range("P26:Q26").UnMerge()
$ws.Range("A25:Q25").Copy()
$ws.Range("A26:Q26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

25..21 | ForEach-Object {
    $src = $_
    $dst = $_ + 1
    $ws.Range("C$dst").Value2 = $ws.Range("C$src").Value2
    $ws.Range("H$dst").Value2 = $ws.Range("H$src").Value2
    Set-TextValue $ws.Range("L$dst") $ws.Range("L$src").Value2
    $ws.Range("N$dst").Value2 = $ws.Range("N$src").Value2
    Set-TextValue $ws.Range("P$dst") $ws.Range("P$src").Value2
    $ws.Range("Q$dst").Value2 = $ws.Range("Q$src").Value2
}

# Row 26 is the item that used to be #19 (the former last row); give it its
# own item number (20) -- column A was otherwise untouched by the loop
# above since item numbers 15-19 in A21:A25 do not change.
$ws.Range("A26").Value2 = 20

# ------------------------------------------------------------------
# 3. Write the new item into row 21 (everything else has been pushed down
#    to make space for it).
# ------------------------------------------------------------------
$ws.Range("C21").Value2 = "THYROXINE 50MCG 100 TAB."
$ws.Range("H21").Value2 = "9:0"
Set-TextValue $ws.Range("L21") "1"
$ws.Range("N21").Value2 = "46.00"
Set-TextValue $ws.Range("P21") "46.0000"
$ws.Range("Q21").Value2 = "1:0"

# ------------------------------------------------------------------
# 4. Update the totals row (now row 27): add the new item's selling price
#    to the previous total.
# ------------------------------------------------------------------
$ws.Range("P27").Value2 = 1553

# ------------------------------------------------------------------
# 5. Refresh the footer timestamp (now row 28).
# ------------------------------------------------------------------
$ws.Range("A28").Value2 = "Monday, 22 September, 2025 12:06 PM"

# ------------------------------------------------------------------
# 6. Fix up merged cells: row 26 changes from the single totals merge
#    (P26:Q26, already unmerged above) to the normal 5-block item-row
#    layout, and a new totals merge is created at row 27.
# ------------------------------------------------------------------
$ws.Range("A26:B26").Merge()
$ws.Range("C26:G26").Merge()
$ws.Range("H26:K26").Merge()
$ws.Range("L26:M26").Merge()
$ws.Range("N26:O26").Merge()
$ws.Range("P27:Q27").Merge()

Write-Output "Inserted THYROXINE 50MCG 100 TAB. as item 15; shifted remaining items; updated totals and timestamp."
